$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet holds a lower-triangular distance matrix with 1-indexed city ids
# in row 1 (B1:AW1, the column headers) and in column A (A2:A49, the row
# headers). Convert both label sequences to zero-indexed city ids by
# decrementing every id by one. The distance values themselves (B2:AW49)
# are left untouched.

# Row 1 headers: B1:AW1 currently hold 1..48 -> becomes 0..47
for ($col = 2; $col -le 49; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $cell.Value2 = $cell.Value2 - 1
}

# Column A headers: A2:A49 currently hold 1..48 -> becomes 0..47
for ($row = 2; $row -le 49; $row++) {
    $cell = $ws.Cells.Item($row, 1)
    $cell.Value2 = $cell.Value2 - 1
}
